# Apply column B ("Numerical"/"Categorical") classification to the variable
# description sheet, extend the used range/autofilter/defined name from
# rows 1:2 to 1:82, and update the worksheet selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> classification label for column B, rows 2 through 82.
$values = @(
    "2:Numerical", "3:Numerical", "4:Categorical", "5:Numerical", "6:Numerical", "7:Categorical", "8:Categorical", "9:Categorical",
    "10:Categorical", "11:Categorical", "12:Categorical", "13:Categorical", "14:Categorical", "15:Categorical", "16:Categorical", "17:Categorical",
    "18:Categorical", "19:Numerical", "20:Numerical", "21:Numerical", "22:Numerical", "23:Categorical", "24:Categorical", "25:Categorical",
    "26:Categorical", "27:Categorical", "28:Numerical", "29:Categorical", "30:Categorical", "31:Categorical", "32:Categorical", "33:Categorical",
    "34:Categorical", "35:Categorical", "36:Numerical", "37:Categorical", "38:Numerical", "39:Numerical", "40:Numerical", "41:Categorical",
    "42:Categorical", "43:Categorical", "44:Categorical", "45:Numerical", "46:Numerical", "47:Numerical", "48:Numerical", "49:Numerical",
    "50:Numerical", "51:Numerical", "52:Numerical", "53:Numerical", "54:Numerical", "55:Categorical", "56:Numerical", "57:Categorical",
    "58:Numerical", "59:Categorical", "60:Categorical", "61:Numerical", "62:Categorical", "63:Numerical", "64:Numerical", "65:Categorical",
    "66:Categorical", "67:Categorical", "68:Numerical", "69:Numerical", "70:Numerical", "71:Numerical", "72:Numerical", "73:Numerical",
    "74:Categorical", "75:Categorical", "76:Categorical", "77:Numerical", "78:Numerical", "79:Numerical", "80:Categorical", "81:Categorical",
    "82:Numerical"
)

foreach ($entry in $values) {
    $parts = $entry.Split(":")
    $row = [int]$parts[0]
    $label = $parts[1]
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $label
    $cell.Style = "Normal"
}

# Extend the AutoFilter to cover the new data range A1:G82.
$ws.AutoFilterMode = $false
$ws.Range("A1:G82").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$82"
    }
}

# Move the active selection to B83, just past the newly written data.
$ws.Range("B83").Select()
